$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 324 - this shifts the existing rows 324:455
# down to 325:456 (and the used range / dimension grows to R456),
# exactly like a weekly re-sort that pushes a new data point to the
# top of this sub-series and bumps every later row down by one.
$ws.Rows.Item(324).Insert()

# Populate the newly inserted row 324 with the new data point.
$ws.Range("A324").Value = 4
$ws.Range("B324").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C324").Value = "Los Lagos"
$ws.Range("D324").Value = 45009
$ws.Range("E324").Value = 10
$ws.Range("F324").Value = 100114014
$ws.Range("G324").Value = "Betarraga"
$ws.Range("H324").Value = "Sin especificar"
$ws.Range("I324").Value = "Primera"
$ws.Range("J324").Value = 1000
$ws.Range("K324").Value = 1100
$ws.Range("L324").Value = 1100
$ws.Range("M324").Value = 1100
$ws.Range("N324").Value = "$/paquete 5 unidades"
$ws.Range("O324").Value = "Provincia de Cautín"
$ws.Range("P324").Value = 220
$ws.Range("Q324").Value = 5
$ws.Range("R324").Value = "Hortaliza"

# Make sure the date cell keeps the sheet's date number format
# (column D elsewhere uses numFmtId 165 "YYYY-MM-DD HH:MM:SS").
$ws.Range("D324").NumberFormat = $ws.Range("D325").NumberFormat
